# Scheduled-runner update: refresh cached Universalis market-board figures
# (currentAveragePrice / NQ / HQ and derived Leve profit columns H:N)
# for the Balmung_Profits workbook, sheet by sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1069.875
$ws.Range("I88").Value = 1224
$ws.Range("J88").Value = 813
$ws.Range("K88").Value = 1224
$ws.Range("L88").Value = 813
$ws.Range("M88").Value = -818
$ws.Range("N88").Value = -1625
$ws.Range("H91").Value = 1069.875
$ws.Range("I91").Value = 1224
$ws.Range("J91").Value = 813
$ws.Range("K91").Value = 1224
$ws.Range("L91").Value = 813
$ws.Range("M91").Value = 180
$ws.Range("N91").Value = -3621
$ws.Range("H92").Value = 2606205
$ws.Range("I92").Value = 2232903.8
$ws.Range("J92").Value = 2843760
$ws.Range("K92").Value = 2232903.8
$ws.Range("L92").Value = 2843760
$ws.Range("M92").Value = -2231655.8
$ws.Range("N92").Value = -2846256
$ws.Range("H112").Value = 335444
$ws.Range("J112").Value = 335444
$ws.Range("L112").Value = 1006332
$ws.Range("N112").Value = -1008548
$ws.Range("H131").Value = 9665.223
$ws.Range("I131").Value = 8873.375
$ws.Range("K131").Value = 26620.125
$ws.Range("M131").Value = -21580.125
$ws.Range("H132").Value = 20893.611
$ws.Range("I132").Value = 27497.2
$ws.Range("J132").Value = 2026.2142
$ws.Range("K132").Value = 82491.60000000001
$ws.Range("L132").Value = 6078.642599999999
$ws.Range("M132").Value = -79961.60000000001
$ws.Range("N132").Value = -11138.6426
$ws.Range("H133").Value = 123105
$ws.Range("J133").Value = 123105
$ws.Range("L133").Value = 123105
$ws.Range("N133").Value = -133225
$ws.Range("H138").Value = 7070.1875
$ws.Range("I138").Value = 11222.934
$ws.Range("J138").Value = 3406
$ws.Range("K138").Value = 33668.802
$ws.Range("L138").Value = 10218
$ws.Range("M138").Value = -28528.802
$ws.Range("N138").Value = -20498
$ws.Range("H139").Value = 159991
$ws.Range("J139").Value = 159991
$ws.Range("L139").Value = 159991
$ws.Range("N139").Value = -170271

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 319555.78
$ws.Range("I32").Value = 373473.2
$ws.Range("K32").Value = 373473.2
$ws.Range("M32").Value = -373186.2
$ws.Range("H45").Value = 103169.2
$ws.Range("I45").Value = 103169.2
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 103169.2
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("M45").Value = -102792.2
$ws.Range("H61").Value = 6145.698
$ws.Range("I61").Value = 5432.0303
$ws.Range("K61").Value = 5432.0303
$ws.Range("M61").Value = -5220.0303
$ws.Range("H104").Value = 56000.332
$ws.Range("J104").Value = 56000.332
$ws.Range("L104").Value = 56000.332
$ws.Range("N104").Value = -62988.332
$ws.Range("H136").Value = 6145.698
$ws.Range("I136").Value = 5432.0303
$ws.Range("K136").Value = 16296.0909
$ws.Range("M136").Value = -13746.0909
$ws.Range("H140").Value = 97658.336
$ws.Range("J140").Value = 97658.336
$ws.Range("L140").Value = 97658.336
$ws.Range("N140").Value = -108018.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1408.7727
$ws.Range("I20").Value = 1265.0834
$ws.Range("J20").Value = 1581.2
$ws.Range("K20").Value = 1265.0834
$ws.Range("L20").Value = 1581.2
$ws.Range("M20").Value = -1018.0834
$ws.Range("N20").Value = -2075.2
$ws.Range("H107").Value = 15078.526
$ws.Range("I107").Value = 20036.46
$ws.Range("K107").Value = 20036.46
$ws.Range("M107").Value = -18116.46

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 6500
$ws.Range("I45").Value = 6500
$ws.Range("K45").Value = 6500
$ws.Range("M45").Value = -5907
$ws.Range("H62").Value = 3010.6667
$ws.Range("I62").Value = 3259.2
$ws.Range("K62").Value = 3259.2
$ws.Range("M62").Value = -2635.2
$ws.Range("H65").Value = 3010.6667
$ws.Range("I65").Value = 3259.2
$ws.Range("K65").Value = 16296
$ws.Range("M65").Value = -13176
$ws.Range("H141").Value = 482968.3
$ws.Range("J141").Value = 431076
$ws.Range("L141").Value = 431076
$ws.Range("N141").Value = -441436

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 401.7143
$ws.Range("J2").Value = 765.625
$ws.Range("L2").Value = 4593.75
$ws.Range("N2").Value = -4819.75
$ws.Range("H25").Value = 2434.8
$ws.Range("J25").Value = 2864.8572
$ws.Range("L25").Value = 8594.571599999999
$ws.Range("N25").Value = -8932.571599999999
$ws.Range("H30").Value = 2434.8
$ws.Range("J30").Value = 2864.8572
$ws.Range("L30").Value = 8594.571599999999
$ws.Range("N30").Value = -8798.571599999999
$ws.Range("H107").Value = 55555990
$ws.Range("J107").Value = 76923390
$ws.Range("L107").Value = 230770170
$ws.Range("N107").Value = -230774010

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.53333
$ws.Range("I2").Value = 58.6
$ws.Range("J2").Value = 153
$ws.Range("K2").Value = 58.6
$ws.Range("L2").Value = 153
$ws.Range("M2").Value = 54.4
$ws.Range("N2").Value = -379
$ws.Range("H3").Value = 5219.8
$ws.Range("I3").Value = 5219.8
$ws.Range("K3").Value = 5219.8
$ws.Range("M3").Value = -5103.8
$ws.Range("H80").Value = 78555.13
$ws.Range("I80").Value = 113032.1
$ws.Range("K80").Value = 113032.1
$ws.Range("M80").Value = -112034.1
$ws.Range("H83").Value = 78555.13
$ws.Range("I83").Value = 113032.1
$ws.Range("K83").Value = 565160.5
$ws.Range("M83").Value = -560168.5
$ws.Range("H102").Value = 38662980
$ws.Range("I102").Value = 41768228
$ws.Range("K102").Value = 41768228
$ws.Range("M102").Value = -41766606
$ws.Range("H122").Value = 5516.737
$ws.Range("I122").Value = 3455.1667
$ws.Range("J122").Value = 9050.857
$ws.Range("K122").Value = 10365.5001
$ws.Range("L122").Value = 27152.571
$ws.Range("M122").Value = -7915.500100000001
$ws.Range("N122").Value = -32052.571
$ws.Range("H126").Value = 2457.5
$ws.Range("I126").Value = 2457.5
$ws.Range("K126").Value = 7372.5
$ws.Range("M126").Value = -4902.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 958.6667
$ws.Range("I16").Value = 447.55554
$ws.Range("J16").Value = 2492
$ws.Range("K16").Value = 447.55554
$ws.Range("L16").Value = 2492
$ws.Range("M16").Value = -277.55554
$ws.Range("N16").Value = -2832
$ws.Range("H22").Value = 3191.2666
$ws.Range("I22").Value = 921
$ws.Range("K22").Value = 921
$ws.Range("M22").Value = -626
$ws.Range("H27").Value = 3191.2666
$ws.Range("I27").Value = 921
$ws.Range("K27").Value = 921
$ws.Range("M27").Value = -814
$ws.Range("H46").Value = 4942.25
$ws.Range("J46").Value = 1627.762
$ws.Range("L46").Value = 1627.762
$ws.Range("N46").Value = -2003.762
$ws.Range("H55").Value = 1198.35
$ws.Range("I55").Value = 2275.5715
$ws.Range("J55").Value = 618.3077
$ws.Range("K55").Value = 2275.5715
$ws.Range("L55").Value = 618.3077
$ws.Range("M55").Value = -2102.5715
$ws.Range("N55").Value = -964.3077
$ws.Range("H61").Value = 4880187
$ws.Range("I61").Value = 6453631
$ws.Range("J61").Value = 2510.3
$ws.Range("K61").Value = 6453631
$ws.Range("L61").Value = 2510.3
$ws.Range("M61").Value = -6453429
$ws.Range("N61").Value = -2914.3
$ws.Range("H68").Value = 3938.4375
$ws.Range("I68").Value = 4053.6
$ws.Range("K68").Value = 4053.6
$ws.Range("M68").Value = -3304.6
$ws.Range("H71").Value = 3938.4375
$ws.Range("I71").Value = 4053.6
$ws.Range("K71").Value = 20268
$ws.Range("M71").Value = -16524
$ws.Range("H93").Value = 1938.4286
$ws.Range("I93").Value = 1094.8334
$ws.Range("K93").Value = 1094.8334
$ws.Range("M93").Value = 153.1666
$ws.Range("H113").Value = 4880187
$ws.Range("I113").Value = 6453631
$ws.Range("J113").Value = 2510.3
$ws.Range("K113").Value = 6453631
$ws.Range("L113").Value = 2510.3
$ws.Range("M113").Value = -6451461
$ws.Range("N113").Value = -6850.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26363.637
$ws.Range("H81").Value = 1746.0769
$ws.Range("I81").Value = 1952
$ws.Range("J81").Value = 1059.6666
$ws.Range("K81").Value = 3904
$ws.Range("L81").Value = 2119.3332
$ws.Range("M81").Value = -2843
$ws.Range("N81").Value = -4241.3332
$ws.Range("H84").Value = 1746.0769
$ws.Range("I84").Value = 1952
$ws.Range("J84").Value = 1059.6666
$ws.Range("K84").Value = 19520
$ws.Range("L84").Value = 10596.666
$ws.Range("M84").Value = -14216
$ws.Range("N84").Value = -21204.666
$ws.Range("H93").Value = 41997.5
$ws.Range("J93").Value = 41997.5
$ws.Range("L93").Value = 41997.5
$ws.Range("N93").Value = -46989.5
$ws.Range("H96").Value = 1279.381
$ws.Range("I96").Value = 791.5
$ws.Range("J96").Value = 1929.8889
$ws.Range("K96").Value = 791.5
$ws.Range("L96").Value = 1929.8889
$ws.Range("M96").Value = 581.5
$ws.Range("N96").Value = -4675.8889
$ws.Range("H104").Value = 31953.143
$ws.Range("J104").Value = 31953.143
$ws.Range("L104").Value = 31953.143
$ws.Range("N104").Value = -38941.143
$ws.Range("H122").Value = 21211.062
$ws.Range("I122").Value = 22492.133
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 67476.399
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -65026.399
$ws.Range("N122").Value = -10885
$ws.Range("H125").Value = 68995.60000000001
$ws.Range("J125").Value = 68995.60000000001
$ws.Range("L125").Value = 68995.60000000001
$ws.Range("N125").Value = -78835.60000000001
